$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new "Label" header in column H, row 1, matching the header style used by B1:G1
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats: copy formatting/style only
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Label"

# Labels: rows 2-11 correspond to 0 (Control) for rows 2-6, 1 (MDD) for rows 7-11
# Same pattern repeats for rows 12-21 (second block, Iterations=200)
$labels = @(0,0,0,0,0,1,1,1,1,1)

for ($i = 0; $i -lt 10; $i++) {
    $row1 = 2 + $i
    $row2 = 12 + $i
    $ws.Cells.Item($row1, 8).Value = $labels[$i]
    $ws.Cells.Item($row2, 8).Value = $labels[$i]
}
